# Melón (Hortaliza, Terminal Hortofrutícola Agro Chillán) — add a new weekly
# report block. Insert 6 new data rows right before the existing row 266
# (pushing the existing 266:356 block down to 272:362), then fill the new
# rows with the new week's Calameño/Tuna x Extra/Primera/Segunda entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 blank rows above the current row 266 — this shifts rows 266:356
# down to 272:362, carrying their formatting (incl. the date-format style on
# column D) along with them.
$ws.Rows("266:271").Insert()

$mercadoId = 7
$mercado   = "Terminal Hortofrutícola Agro Chillán"
$region    = "Ñuble"
$fecha     = 44992
$codreg    = 16
$categoriaId = 100112027
$categoria   = "Melón"
$unidad      = "`$/unidad"
$kgOUnidades = 1
$clasificacion = "Hortaliza"

# r, Variedad, Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm, Origen, PrecioKg
$rows = @(
    @(266, "Calameño", "Extra",   1000, 800, 800, 800, "Región del Maule", 800),
    @(267, "Calameño", "Primera", 1000, 600, 600, 600, "Región del Maule", 600),
    @(268, "Calameño", "Segunda", 1000, 400, 400, 400, "Región del Maule", 400),
    @(269, "Tuna",     "Extra",   1000, 800, 800, 800, "Región del Maule", 800),
    @(270, "Tuna",     "Primera", 1000, 600, 600, 600, "Región del Maule", 600),
    @(271, "Tuna",     "Segunda", 1000, 400, 400, 400, "Región del Maule", 400)
)

foreach ($r in $rows) {
    $rowNum      = $r[0]
    $variedad    = $r[1]
    $calidad     = $r[2]
    $volumen     = $r[3]
    $precioMin   = $r[4]
    $precioMax   = $r[5]
    $precioProm  = $r[6]
    $origen      = $r[7]
    $precioKg    = $r[8]

    $ws.Cells.Item($rowNum, 1).Value  = $mercadoId
    $ws.Cells.Item($rowNum, 2).Value  = $mercado
    $ws.Cells.Item($rowNum, 3).Value  = $region
    $ws.Cells.Item($rowNum, 4).Value  = $fecha
    $ws.Cells.Item($rowNum, 5).Value  = $codreg
    $ws.Cells.Item($rowNum, 6).Value  = $categoriaId
    $ws.Cells.Item($rowNum, 7).Value  = $categoria
    $ws.Cells.Item($rowNum, 8).Value  = $variedad
    $ws.Cells.Item($rowNum, 9).Value  = $calidad
    $ws.Cells.Item($rowNum, 10).Value = $volumen
    $ws.Cells.Item($rowNum, 11).Value = $precioMin
    $ws.Cells.Item($rowNum, 12).Value = $precioMax
    $ws.Cells.Item($rowNum, 13).Value = $precioProm
    $ws.Cells.Item($rowNum, 14).Value = $unidad
    $ws.Cells.Item($rowNum, 15).Value = $origen
    $ws.Cells.Item($rowNum, 16).Value = $precioKg
    $ws.Cells.Item($rowNum, 17).Value = $kgOUnidades
    $ws.Cells.Item($rowNum, 18).Value = $clasificacion
}
